# Correção nos dados e início da análise PNAD 2009
#
# Changes applied:
#  1. Row 2 header cells B2 and F2 were placeholder pandas column names
#     ("unnamed: 1_level_1" / "unnamed: 5_level_1"); they are corrected to
#     "total" (matching the other "total" header already in C2).
#  2. Two stray section-header rows that only contained a label in column A
#     and no data ("situação do domicílio" at row 5 and "grandes regiões e
#     unidades da federação" at row 8) are removed entirely, so the data
#     rows below shift up to close the gaps.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("dados")

# 1) Fix the placeholder header labels in row 2.
$ws.Range("B2").Value = "total"
$ws.Range("F2").Value = "total"

# 2) Delete the two empty label-only rows (higher row number first so the
#    second deletion's row index is unaffected by the first).
$ws.Rows.Item(8).Delete()
$ws.Rows.Item(5).Delete()
